# Insert a new row at position 113 (shifts existing rows 113..223 down to 114..224)
# and populate it with the new weekly price record (2021-10-15 / serial 44484).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(113).Insert()

$ws.Cells.Item(113, 1).Value = 9
$ws.Cells.Item(113, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(113, 3).Value = "Metropolitana"
$ws.Cells.Item(113, 4).Value = 44484
$ws.Cells.Item(113, 5).Value = 13
$ws.Cells.Item(113, 6).Value = 100112044
$ws.Cells.Item(113, 7).Value = "Perejil"
$ws.Cells.Item(113, 8).Value = "Sin especificar"
$ws.Cells.Item(113, 9).Value = "Primera"
$ws.Cells.Item(113, 10).Value = 97
$ws.Cells.Item(113, 11).Value = 8000
$ws.Cells.Item(113, 12).Value = 10000
$ws.Cells.Item(113, 13).Value = 8990
$ws.Cells.Item(113, 14).Value = "$/docena de atados"
$ws.Cells.Item(113, 15).Value = "Región Metropolitana"
$ws.Cells.Item(113, 16).Value = 2997
$ws.Cells.Item(113, 17).Value = 3
$ws.Cells.Item(113, 18).Value = "Hortaliza"
